$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 573; this shifts existing rows 573..621 down to 574..622
$ws.Rows.Item(573).Insert()

# Populate the new row 573 with the new data entry.
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg, F Categoria ID,
# G Categoria, H Variedad, I Calidad, J Volumen, K Precio minimo, L Precio maximo,
# M Precio promedio ponderado, N Unidad de comercializacion, O Origen,
# P Precio $/Kg, Q Kg o Unidades, R Clasificacion
$ws.Range("A573").Value = 4
$ws.Range("B573").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C573").Value = "Los Lagos"
$ws.Range("D573").Value = 44769
$ws.Range("E573").Value = 10
$ws.Range("F573").Value = 100112004
$ws.Range("G573").Value = "Cebolla"
$ws.Range("H573").Value = "Sin especificar"
$ws.Range("I573").Value = "1a (guarda)"
$ws.Range("J573").Value = 150
$ws.Range("K573").Value = 10000
$ws.Range("L573").Value = 10000
$ws.Range("M573").Value = 10000
$ws.Range("N573").Value = "`$/malla 18 kilos"
$ws.Range("O573").Value = "Región de O'Higgins"
$ws.Range("P573").Value = 556
$ws.Range("Q573").Value = 18
$ws.Range("R573").Value = "Hortaliza"
